$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lapas1")

# Add a new row of data: B15 = text, C15 = number
$ws.Range("B15").Value = "Paveiksleliu api taisymas ir frontendas"
$ws.Range("C15").Value = 3
